# Auto-generated edit script: re-sort/update Serie C Group B 2023-2024 match rows
# and append the newly scraped Vis Pesaro vs Torres fixture as row 137.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows whose match data moved position (columns A-E, i.e. Indice/pais/torneio/temporada/data_partida, are unchanged) ---
# Row 93: Ancona 3-0 Arezzo
$ws.Range("F93").Value = "Ancona"
$ws.Range("G93").Value = 3
$ws.Range("H93").Value = "Arezzo"
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 2.16
$ws.Range("K93").Value = "24/10/2023 21:12"
$ws.Range("L93").Value = 1.94
$ws.Range("M93").Value = "26/10/2023 18:22"
$ws.Range("N93").Value = 2.9
$ws.Range("O93").Value = "24/10/2023 21:12"
$ws.Range("P93").Value = 3.4
$ws.Range("Q93").Value = "26/10/2023 18:22"
$ws.Range("R93").Value = 3.41
$ws.Range("S93").Value = "24/10/2023 21:12"
$ws.Range("T93").Value = 4
$ws.Range("U93").Value = "26/10/2023 18:22"
$ws.Range("V93").Value = "https://www.betexplorer.com/football/italy/serie-c-group-b/ancona-arezzo/lp1QsBv0/"

# Row 94: Carrarese 1-1 Pineto
$ws.Range("F94").Value = "Carrarese"
$ws.Range("G94").Value = 1
$ws.Range("H94").Value = "Pineto"
$ws.Range("I94").Value = 1
$ws.Range("J94").Value = 1.67
$ws.Range("K94").Value = "24/10/2023 21:12"
$ws.Range("L94").Value = 1.52
$ws.Range("M94").Value = "26/10/2023 18:21"
$ws.Range("N94").Value = 3.29
$ws.Range("O94").Value = "24/10/2023 21:12"
$ws.Range("P94").Value = 3.84
$ws.Range("Q94").Value = "26/10/2023 18:21"
$ws.Range("R94").Value = 5.03
$ws.Range("S94").Value = "24/10/2023 21:12"
$ws.Range("T94").Value = 7.15
$ws.Range("U94").Value = "26/10/2023 18:21"
$ws.Range("V94").Value = "https://www.betexplorer.com/football/italy/serie-c-group-b/carrarese-pineto/W8jVtVg6/"

# Row 95: Pontedera 0-0 Vis Pesaro
$ws.Range("F95").Value = "Pontedera"
$ws.Range("G95").Value = 0
$ws.Range("H95").Value = "Vis Pesaro"
$ws.Range("I95").Value = 0
$ws.Range("J95").Value = 1.84
$ws.Range("K95").Value = "24/10/2023 21:12"
$ws.Range("L95").Value = 1.82
$ws.Range("M95").Value = "26/10/2023 18:02"
$ws.Range("N95").Value = 3.08
$ws.Range("O95").Value = "24/10/2023 21:12"
$ws.Range("P95").Value = 3.34
$ws.Range("Q95").Value = "26/10/2023 18:02"
$ws.Range("R95").Value = 4.28
$ws.Range("S95").Value = "24/10/2023 21:12"
$ws.Range("T95").Value = 4.8
$ws.Range("U95").Value = "26/10/2023 18:02"
$ws.Range("V95").Value = "https://www.betexplorer.com/football/italy/serie-c-group-b/us-pontedera-vis-pesaro/McJNnLsP/"

# Row 97: Spal 1-0 Sestri Levante
$ws.Range("F97").Value = "Spal"
$ws.Range("G97").Value = 1
$ws.Range("H97").Value = "Sestri Levante"
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 1.76
$ws.Range("K97").Value = "24/10/2023 21:12"
$ws.Range("L97").Value = 1.85
$ws.Range("M97").Value = "26/10/2023 18:29"
$ws.Range("N97").Value = 3.21
$ws.Range("O97").Value = "24/10/2023 21:12"
$ws.Range("P97").Value = 3.03
$ws.Range("Q97").Value = "26/10/2023 18:29"
$ws.Range("R97").Value = 4.72
$ws.Range("S97").Value = "24/10/2023 21:12"
$ws.Range("T97").Value = 5.29
$ws.Range("U97").Value = "26/10/2023 18:29"
$ws.Range("V97").Value = "https://www.betexplorer.com/football/italy/serie-c-group-b/spal-sestri-levante/WW8ctJKh/"

# Row 99: Torres 1-1 Spal
$ws.Range("F99").Value = "Torres"
$ws.Range("G99").Value = 1
$ws.Range("H99").Value = "Spal"
$ws.Range("I99").Value = 1
$ws.Range("J99").Value = 1.79
$ws.Range("K99").Value = "26/10/2023 22:12"
$ws.Range("L99").Value = 2
$ws.Range("M99").Value = "29/10/2023 13:53"
$ws.Range("N99").Value = 3.12
$ws.Range("O99").Value = "26/10/2023 22:12"
$ws.Range("P99").Value = 3.18
$ws.Range("Q99").Value = "29/10/2023 13:53"
$ws.Range("R99").Value = 4.5
$ws.Range("S99").Value = "26/10/2023 22:12"
$ws.Range("T99").Value = 4.12
$ws.Range("U99").Value = "29/10/2023 13:54"
$ws.Range("V99").Value = "https://www.betexplorer.com/football/italy/serie-c-group-b/sassari-torres-spal/Q3uZuAQ9/"

# Row 100: Juventus U23 3-1 Olbia
$ws.Range("F100").Value = "Juventus U23"
$ws.Range("G100").Value = 3
$ws.Range("H100").Value = "Olbia"
$ws.Range("I100").Value = 1
$ws.Range("J100").Value = 2.08
$ws.Range("K100").Value = "26/10/2023 22:12"
$ws.Range("L100").Value = 2.06
$ws.Range("M100").Value = "29/10/2023 13:33"
$ws.Range("N100").Value = 2.94
$ws.Range("O100").Value = "26/10/2023 22:12"
$ws.Range("P100").Value = 3.17
$ws.Range("Q100").Value = "29/10/2023 13:50"
$ws.Range("R100").Value = 3.57
$ws.Range("S100").Value = "26/10/2023 22:12"
$ws.Range("T100").Value = 3.6
$ws.Range("U100").Value = "29/10/2023 13:33"
$ws.Range("V100").Value = "https://www.betexplorer.com/football/italy/serie-c-group-b/juventus-olbia/Q7ezL8lp/"

# Row 104: Perugia 2-1 Entella
$ws.Range("F104").Value = "Perugia"
$ws.Range("G104").Value = 2
$ws.Range("H104").Value = "Entella"
$ws.Range("I104").Value = 1
$ws.Range("J104").Value = 2.12
$ws.Range("K104").Value = "26/10/2023 22:12"
$ws.Range("L104").Value = 1.92
$ws.Range("M104").Value = "30/10/2023 20:36"
$ws.Range("N104").Value = 3.01
$ws.Range("O104").Value = "26/10/2023 22:12"
$ws.Range("P104").Value = 3.28
$ws.Range("Q104").Value = "30/10/2023 20:36"
$ws.Range("R104").Value = 3.36
$ws.Range("S104").Value = "26/10/2023 22:12"
$ws.Range("T104").Value = 4.34
$ws.Range("U104").Value = "30/10/2023 20:36"
$ws.Range("V104").Value = "https://www.betexplorer.com/football/italy/serie-c-group-b/perugia-entella/hEmMrCti/"

# Row 105: Cesena 2-1 Carrarese
$ws.Range("F105").Value = "Cesena"
$ws.Range("G105").Value = 2
$ws.Range("H105").Value = "Carrarese"
$ws.Range("I105").Value = 1
$ws.Range("J105").Value = 1.84
$ws.Range("K105").Value = "26/10/2023 21:12"
$ws.Range("L105").Value = 1.65
$ws.Range("M105").Value = "30/10/2023 20:42"
$ws.Range("N105").Value = 3.08
$ws.Range("O105").Value = "26/10/2023 21:12"
$ws.Range("P105").Value = 3.44
$ws.Range("Q105").Value = "30/10/2023 20:42"
$ws.Range("R105").Value = 4.28
$ws.Range("S105").Value = "26/10/2023 21:12"
$ws.Range("T105").Value = 6.35
$ws.Range("U105").Value = "30/10/2023 20:42"
$ws.Range("V105").Value = "https://www.betexplorer.com/football/italy/serie-c-group-b/cesena-carrarese/W6w2wH4B/"

# Row 106: Arezzo 1-0 Gubbio
$ws.Range("F106").Value = "Arezzo"
$ws.Range("G106").Value = 1
$ws.Range("H106").Value = "Gubbio"
$ws.Range("I106").Value = 0
$ws.Range("J106").Value = 2.34
$ws.Range("K106").Value = "26/10/2023 22:12"
$ws.Range("L106").Value = 3.4
$ws.Range("M106").Value = "30/10/2023 20:42"
$ws.Range("N106").Value = 2.85
$ws.Range("O106").Value = "26/10/2023 22:12"
$ws.Range("P106").Value = 2.96
$ws.Range("Q106").Value = "30/10/2023 20:42"
$ws.Range("R106").Value = 3.09
$ws.Range("S106").Value = "26/10/2023 22:12"
$ws.Range("T106").Value = 2.35
$ws.Range("U106").Value = "30/10/2023 20:42"
$ws.Range("V106").Value = "https://www.betexplorer.com/football/italy/serie-c-group-b/arezzo-gubbio/xpvbvck5/"

# Row 107: Vis Pesaro 1-1 Pineto
$ws.Range("F107").Value = "Vis Pesaro"
$ws.Range("G107").Value = 1
$ws.Range("H107").Value = "Pineto"
$ws.Range("I107").Value = 1
$ws.Range("J107").Value = 2.21
$ws.Range("K107").Value = "26/10/2023 22:12"
$ws.Range("L107").Value = 2.65
$ws.Range("M107").Value = "30/10/2023 20:36"
$ws.Range("N107").Value = 2.88
$ws.Range("O107").Value = "26/10/2023 22:12"
$ws.Range("P107").Value = 2.76
$ws.Range("Q107").Value = "30/10/2023 20:36"
$ws.Range("R107").Value = 3.31
$ws.Range("S107").Value = "26/10/2023 22:12"
$ws.Range("T107").Value = 3.16
$ws.Range("U107").Value = "30/10/2023 20:36"
$ws.Range("V107").Value = "https://www.betexplorer.com/football/italy/serie-c-group-b/vis-pesaro-pineto/KQvwuUuG/"

# Row 135: Gubbio 5-2 Sestri Levante
$ws.Range("F135").Value = "Gubbio"
$ws.Range("G135").Value = 5
$ws.Range("H135").Value = "Sestri Levante"
$ws.Range("I135").Value = 2
$ws.Range("J135").Value = 1.53
$ws.Range("K135").Value = "16/11/2023 18:12"
$ws.Range("L135").Value = 1.6
$ws.Range("M135").Value = "19/11/2023 16:11"
$ws.Range("N135").Value = 3.75
$ws.Range("O135").Value = "16/11/2023 18:12"
$ws.Range("P135").Value = 3.61
$ws.Range("Q135").Value = "19/11/2023 16:11"
$ws.Range("R135").Value = 5.85
$ws.Range("S135").Value = "16/11/2023 18:12"
$ws.Range("T135").Value = 6.56
$ws.Range("U135").Value = "19/11/2023 16:11"
$ws.Range("V135").Value = "https://www.betexplorer.com/football/italy/serie-c-group-b/gubbio-sestri-levante/EVApdMG7/"

# Row 136: Carrarese 1-0 Spal
$ws.Range("F136").Value = "Carrarese"
$ws.Range("G136").Value = 1
$ws.Range("H136").Value = "Spal"
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 1.89
$ws.Range("K136").Value = "16/11/2023 09:12"
$ws.Range("L136").Value = 1.66
$ws.Range("M136").Value = "19/11/2023 16:12"
$ws.Range("N136").Value = 3.15
$ws.Range("O136").Value = "16/11/2023 09:12"
$ws.Range("P136").Value = 3.43
$ws.Range("Q136").Value = "19/11/2023 16:12"
$ws.Range("R136").Value = 3.89
$ws.Range("S136").Value = "16/11/2023 09:12"
$ws.Range("T136").Value = 6.12
$ws.Range("U136").Value = "19/11/2023 16:12"
$ws.Range("V136").Value = "https://www.betexplorer.com/football/italy/serie-c-group-b/carrarese-spal/SE8YbOWl/"

# --- Append brand-new row 137 (copy cell style from row 136 for the Indice/data_partida columns) ---
$ws.Range("A136").Copy() | Out-Null
$ws.Range("A137").PasteSpecial(-4122) | Out-Null
$ws.Range("E136").Copy() | Out-Null
$ws.Range("E137").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws.Range("A137").Value = 136
$ws.Range("B137").Value = "italy"
$ws.Range("C137").Value = "serie-c-group-b"
$ws.Range("D137").Value = "2023-2024"
$ws.Range("E137").Value = 45250.77083333334
$ws.Range("F137").Value = "Vis Pesaro"
$ws.Range("G137").Value = 1
$ws.Range("H137").Value = "Torres"
$ws.Range("I137").Value = 2
$ws.Range("J137").Value = 3.49
$ws.Range("K137").Value = "16/11/2023 09:12"
$ws.Range("L137").Value = 4.27
$ws.Range("M137").Value = "20/11/2023 18:25"
$ws.Range("N137").Value = 3.04
$ws.Range("O137").Value = "16/11/2023 09:12"
$ws.Range("P137").Value = 3.04
$ws.Range("Q137").Value = "20/11/2023 18:25"
$ws.Range("R137").Value = 2.06
$ws.Range("S137").Value = "16/11/2023 09:12"
$ws.Range("T137").Value = 2.02
$ws.Range("U137").Value = "20/11/2023 18:25"
$ws.Range("V137").Value = "https://www.betexplorer.com/football/italy/serie-c-group-b/vis-pesaro-sassari-torres/6Le0DNor/"
